$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# --- Metadata sheet -------------------------------------------------------
# Insert a new "Jurisdiction" property row before "Description" (currently
# row 11), pushing Description/Purpose/Copyright/Immutable down by one.
$wsMeta.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows.
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"

# Refresh the build "Date" value (row 8, column B).
$wsMeta.Range("B8").Value = "2024-07-01T07:50:29+00:00"
